$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.04011237248778343
$ws.Range("C2").Value = 0.017987513914704323
$ws.Range("D2").Value = 0.011631176806986332
$ws.Range("E2").Value = 0.014233378693461418
$ws.Range("F2").Value = 0.00023576710373163223
$ws.Range("G2").Value = 0.002413670066744089
$ws.Range("J2").Value = 0.12776559591293335
$ws.Range("K2").Value = 1.4763305187225342
